$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 557.8570999999999
$ws.Range("I28").Value = 550.8333
$ws.Range("K28").Value = 550.8333
$ws.Range("M28").Value = -65.83330000000001
# Row 64
$ws.Range("H64").Value = 4905.222
$ws.Range("I64").Value = 4993.4375
$ws.Range("J64").Value = 4199.5
$ws.Range("K64").Value = 4993.4375
$ws.Range("L64").Value = 4199.5
$ws.Range("M64").Value = -4745.4375
$ws.Range("N64").Value = -4695.5
# Row 67
$ws.Range("H67").Value = 4905.222
$ws.Range("I67").Value = 4993.4375
$ws.Range("J67").Value = 4199.5
$ws.Range("K67").Value = 4993.4375
$ws.Range("L67").Value = 4199.5
$ws.Range("M67").Value = -4135.4375
$ws.Range("N67").Value = -5915.5
# Row 107
$ws.Range("H107").Value = 97.5
$ws.Range("I107").Value = 97.5
$ws.Range("K107").Value = 97.5
$ws.Range("M107").Value = 1822.5
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
# Row 132
$ws.Range("H132").Value = 2185.1155
$ws.Range("I132").Value = 2072.52
$ws.Range("K132").Value = 6217.559999999999
$ws.Range("M132").Value = -3687.559999999999
# Row 137
$ws.Range("H137").Value = 3951
$ws.Range("I137").Value = 3932.3333
$ws.Range("K137").Value = 11796.9999
$ws.Range("M137").Value = -9246.999899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3425.25
$ws.Range("J2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("N2").Value = -3726
# Row 8
$ws.Range("H8").Value = 433.33334
$ws.Range("I8").Value = 550
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 550
$ws.Range("L8").Value = 200
$ws.Range("M8").Value = -406
$ws.Range("N8").Value = -488
# Row 32
$ws.Range("H32").Value = 4450.1816
$ws.Range("I32").Value = 4145.2
$ws.Range("K32").Value = 4145.2
$ws.Range("M32").Value = -3858.2
# Row 74
$ws.Range("H74").Value = 2374.5
$ws.Range("I74").Value = 1599.2
$ws.Range("K74").Value = 1599.2
$ws.Range("M74").Value = -725.2
# Row 77
$ws.Range("H77").Value = 2374.5
$ws.Range("I77").Value = 1599.2
$ws.Range("K77").Value = 7996
$ws.Range("M77").Value = -3628
# Row 116
$ws.Range("H116").Value = 3425.25
$ws.Range("J116").Value = 3500
$ws.Range("L116").Value = 3500
$ws.Range("N116").Value = -8088

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3425.25
$ws.Range("J3").Value = 3500
$ws.Range("L3").Value = 3500
$ws.Range("N3").Value = -3728
# Row 10
$ws.Range("H10").Value = 257
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
# Row 86
$ws.Range("H86").Value = 9997.5
$ws.Range("J86").Value = 9995
$ws.Range("L86").Value = 9995
$ws.Range("N86").Value = -12241
# Row 89
$ws.Range("H89").Value = 9997.5
$ws.Range("J89").Value = 9995
$ws.Range("L89").Value = 49975
$ws.Range("N89").Value = -61207
# Row 94
$ws.Range("H94").Value = 1708.0952
$ws.Range("I94").Value = 1630.5555
$ws.Range("K94").Value = 1630.5555
$ws.Range("M94").Value = -1179.5555
# Row 105
$ws.Range("H105").Value = 3172.8
$ws.Range("I105").Value = 3029.5
$ws.Range("J105").Value = 3459.4
$ws.Range("K105").Value = 3029.5
$ws.Range("L105").Value = 3459.4
$ws.Range("M105").Value = -1282.5
$ws.Range("N105").Value = -6953.4
# Row 134
$ws.Range("H134").Value = 3496.5
$ws.Range("I134").Value = 3496.5
$ws.Range("K134").Value = 10489.5
$ws.Range("M134").Value = -7954.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1046.8334
$ws.Range("I16").Value = 1087.909
$ws.Range("K16").Value = 1087.909
$ws.Range("M16").Value = -800.9090000000001
# Row 113
$ws.Range("H113").Value = 1046.8334
$ws.Range("I113").Value = 1087.909
$ws.Range("K113").Value = 1087.909
$ws.Range("M113").Value = 1082.091
# Row 122
$ws.Range("H122").Value = 837.25
$ws.Range("I122").Value = 837.25
$ws.Range("K122").Value = 2511.75
$ws.Range("M122").Value = -61.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2999
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 113
$ws.Range("H113").Value = 2605.5
$ws.Range("I113").Value = 2605.5
$ws.Range("K113").Value = 2605.5
$ws.Range("M113").Value = -435.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 55736.4
$ws.Range("I7").Value = 55736.4
$ws.Range("K7").Value = 55736.4
$ws.Range("M7").Value = -55624.4
# Row 16
$ws.Range("H16").Value = 2287.125
$ws.Range("I16").Value = 2287.125
$ws.Range("K16").Value = 2287.125
$ws.Range("M16").Value = -2117.125
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 68
$ws.Range("H68").Value = 52497.5
$ws.Range("J68").Value = 101499.5
$ws.Range("L68").Value = 101499.5
$ws.Range("N68").Value = -102997.5
# Row 71
$ws.Range("H71").Value = 52497.5
$ws.Range("J71").Value = 101499.5
$ws.Range("L71").Value = 507497.5
$ws.Range("N71").Value = -514985.5
# Row 115
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350
# Row 122
$ws.Range("H122").Value = 3503.125
$ws.Range("I122").Value = 3360
$ws.Range("J122").Value = 4505
$ws.Range("K122").Value = 10080
$ws.Range("L122").Value = 13515
$ws.Range("M122").Value = -7630
$ws.Range("N122").Value = -18415
# Row 126
$ws.Range("H126").Value = 55736.4
$ws.Range("I126").Value = 55736.4
$ws.Range("K126").Value = 167209.2
$ws.Range("M126").Value = -164739.2
# Row 132
$ws.Range("H132").Value = 4701.353
$ws.Range("I132").Value = 1981.1111
$ws.Range("K132").Value = 5943.3333
$ws.Range("M132").Value = -3413.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3055.7144
$ws.Range("I81").Value = 3055.7144
$ws.Range("K81").Value = 6111.4288
$ws.Range("M81").Value = -5050.4288
# Row 84
$ws.Range("H84").Value = 3055.7144
$ws.Range("I84").Value = 3055.7144
$ws.Range("K84").Value = 30557.144
$ws.Range("M84").Value = -25253.144
# Row 107
$ws.Range("H107").Value = 481.91666
$ws.Range("I107").Value = 219.8
$ws.Range("K107").Value = 659.4000000000001
$ws.Range("M107").Value = 1260.6
# Row 122
$ws.Range("H122").Value = 2735.3333
$ws.Range("I122").Value = 2735.3333
$ws.Range("K122").Value = 8205.999899999999
$ws.Range("M122").Value = -5755.999899999999
